$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sheet1 = $wb.Worksheets.Item(1)

# Row 3
$c = $ws.Cells.Item(3,1); $c.Value = 112105307; $c.Style = "Normal"
$c = $ws.Cells.Item(3,2); $c.Value = 88966; $c.Style = "Normal"
$c = $ws.Cells.Item(3,3); $c.Value = "'Ovaliderad"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,4); $c.Value = "'NT"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,5); $c.Value = 5754; $c.Style = "Normal"
$c = $ws.Cells.Item(3,6); $c.Value = "'Gultoppig fingersvamp"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,7); $c.Value = "'Ramaria testaceoflava"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,8); $c.Value = "'(Bres.) Corner"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,9); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,10); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,11); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,14); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,16); $c.Value = "'Landverktjärnen (Landverktjärnen), Jmt"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,17); $c.Value = 446544.2690787801; $c.Style = "Normal"
$c = $ws.Cells.Item(3,18); $c.Value = 7032737.655252043; $c.Style = "Normal"
$c = $ws.Cells.Item(3,19); $c.Value = 10; $c.Style = "Normal"
$c = $ws.Cells.Item(3,20); $c.Value = "'Jämtland"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,21); $c.Value = "'Krokom"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,22); $c.Value = "'Jämtland"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,23); $c.Value = "'Alsen"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,25); $c.Value = "'2023-09-15"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,26); $c.Value = "'00:00"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,27); $c.Value = "'2023-09-15"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,28); $c.Value = "'00:00"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,30); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(3,31); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(3,32); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,33); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(3,46); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,49); $c.Value = "'Rashid Kadhim"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,50); $c.Value = "'Rashid Kadhim"; $c.Style = "Normal"
$c = $ws.Cells.Item(3,51); $c.Value = "'"; $c.Style = "Normal"

# Row 4
$c = $ws.Cells.Item(4,1); $c.Value = 112105366; $c.Style = "Normal"
$c = $ws.Cells.Item(4,2); $c.Value = 84808; $c.Style = "Normal"
$c = $ws.Cells.Item(4,3); $c.Value = "'Ovaliderad"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,4); $c.Value = "'LC"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,5); $c.Value = 236435; $c.Style = "Normal"
$c = $ws.Cells.Item(4,6); $c.Value = "'Droppklibbskivling"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,7); $c.Value = "'Limacella guttata"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,8); $c.Value = "'(Pers. : Fr.) Konrad & Maubl."; $c.Style = "Normal"
$c = $ws.Cells.Item(4,9); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,10); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,11); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,14); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,16); $c.Value = "'Landverktjärnen (Landverktjärnen), Jmt"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,17); $c.Value = 446547.3108286796; $c.Style = "Normal"
$c = $ws.Cells.Item(4,18); $c.Value = 7032731.78096032; $c.Style = "Normal"
$c = $ws.Cells.Item(4,19); $c.Value = 10; $c.Style = "Normal"
$c = $ws.Cells.Item(4,20); $c.Value = "'Jämtland"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,21); $c.Value = "'Krokom"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,22); $c.Value = "'Jämtland"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,23); $c.Value = "'Alsen"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,25); $c.Value = "'2023-09-15"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,26); $c.Value = "'00:00"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,27); $c.Value = "'2023-09-15"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,28); $c.Value = "'00:00"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,30); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(4,31); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(4,32); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,33); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(4,46); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,49); $c.Value = "'Rashid Kadhim"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,50); $c.Value = "'Rashid Kadhim"; $c.Style = "Normal"
$c = $ws.Cells.Item(4,51); $c.Value = "'"; $c.Style = "Normal"

# Row 5
$c = $ws.Cells.Item(5,1); $c.Value = 112105682; $c.Style = "Normal"
$c = $ws.Cells.Item(5,2); $c.Value = 88956; $c.Style = "Normal"
$c = $ws.Cells.Item(5,3); $c.Value = "'Ovaliderad"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,4); $c.Value = "'VU"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,5); $c.Value = 5747; $c.Style = "Normal"
$c = $ws.Cells.Item(5,6); $c.Value = "'Läderdoftande fingersvamp"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,7); $c.Value = "'Ramaria safraniolens"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,8); $c.Value = "'Christian"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,9); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,10); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,11); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,14); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,16); $c.Value = "'Svensbergsbäcken (Svensbergsbäcken), Jmt"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,17); $c.Value = 446627.1477181673; $c.Style = "Normal"
$c = $ws.Cells.Item(5,18); $c.Value = 7032919.480234488; $c.Style = "Normal"
$c = $ws.Cells.Item(5,19); $c.Value = 10; $c.Style = "Normal"
$c = $ws.Cells.Item(5,20); $c.Value = "'Jämtland"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,21); $c.Value = "'Krokom"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,22); $c.Value = "'Jämtland"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,23); $c.Value = "'Alsen"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,25); $c.Value = "'2023-09-15"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,26); $c.Value = "'00:00"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,27); $c.Value = "'2023-09-15"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,28); $c.Value = "'00:00"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,30); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(5,31); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(5,32); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,33); $c.Value = $false; $c.Style = "Normal"
$c = $ws.Cells.Item(5,46); $c.Value = "'"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,49); $c.Value = "'Rashid Kadhim"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,50); $c.Value = "'Rashid Kadhim"; $c.Style = "Normal"
$c = $ws.Cells.Item(5,51); $c.Value = "'"; $c.Style = "Normal"

